$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values as plain text
# (e.g. "11.10", "0.0500") so trailing zeros / thousands-dot formatting
# survive. Assigning a numeric-looking string straight to .Value lets
# Excel auto-detect it as a real number and silently normalize it
# (e.g. "11.10" -> 11.1). Force those specific cells to Text format
# first so the new value is stored verbatim as a string, matching the
# existing data convention in this sheet.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '72.292.12'
$ws.Range("E2").Value = '  +4.25%  '
$ws.Range("D3").Value = '4.044.30'
$ws.Range("E3").Value = '  +3.64%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '519.96'
$ws.Range("E5").Value = '  -1.84%  '
$ws.Range("D6").Value = '147.41'
$ws.Range("E6").Value = '  +1.59%  '
$ws.Range("D7").Value = '0.724'
$ws.Range("E7").Value = '  +18.29%  '
$ws.Range("D8").Value = '4.036.10'
$ws.Range("E8").Value = '  +3.61%  '
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").Value = '0.777'
$ws.Range("E10").Value = '  +8.23%  '
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").Value = '0.0000328'
$ws.Range("E12").Value = '  -2.79%  '
$ws.Range("D13").Value = '47.73'
$ws.Range("E13").Value = '  +13.08%  '
$ws.Range("D14").Value = '11.10'
$ws.Range("E14").Value = '  +8.14%  '
$ws.Range("D15").Value = '4.684.96'
$ws.Range("E15").Value = '  +3.64%  '
$ws.Range("D16").Value = '4.057.79'
$ws.Range("E16").Value = '  +4.32%  '
$ws.Range("D17").Value = '21.18'
$ws.Range("E17").Value = '  +6.74%  '
$ws.Range("D18").Value = '14.12'
$ws.Range("E18").Value = '  +0.87%  '
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("E20").Value = '  -0.79%  '
$ws.Range("D21").Value = '72.178.52'
$ws.Range("E21").Value = '  +4.23%  '
$ws.Range("D22").Value = '443.73'
$ws.Range("E22").Value = '  +4.31%  '
$ws.Range("D23").Value = '104.87'
$ws.Range("E23").Value = '  +18.85%  '
$ws.Range("D24").Value = '3.55'
$ws.Range("E24").Value = '  +4.61%  '
$ws.Range("D25").Value = '14.88'
$ws.Range("E25").Value = '  +5.01%  '
$ws.Range("D26").Value = '4.01'
$ws.Range("E26").Value = '  -0.67%  '
$ws.Range("D27").Value = '11.46'
$ws.Range("E27").Value = '  +0.23%  '
$ws.Range("D28").Value = '11.01'
$ws.Range("E28").Value = '  +3.90%  '
$ws.Range("D29").Value = '37.72'
$ws.Range("E29").Value = '  +3.57%  '
$ws.Range("D30").Value = '5.80'
$ws.Range("E30").Value = '  +2.14%  '
$ws.Range("D31").Value = '3.23'
$ws.Range("E31").Value = '  +15.05%  '
$ws.Range("D32").Value = '13.68'
$ws.Range("E32").Value = '  +3.68%  '
$ws.Range("E33").Value = '  +3.31%  '
$ws.Range("D34").Value = '681.95'
$ws.Range("E34").Value = '  -1.27%  '
$ws.Range("D35").Value = '6.82'
$ws.Range("E35").Value = '  +14.83%  '
$ws.Range("D36").Value = '66.87'
$ws.Range("D37").Value = '42.63'
$ws.Range("E37").Value = '  +6.45%  '
$ws.Range("D38").Value = '0.0₃0865'
$ws.Range("E38").Value = '  -1.75%  '
$ws.Range("D39").Value = '0.428'
$ws.Range("E39").Value = '  -1.87%  '
$ws.Range("D40").Value = '3.54'
$ws.Range("E40").Value = '  +7.54%  '
$ws.Range("D41").Value = '0.151'
$ws.Range("E41").Value = '  +0.66%  '
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("D43").Value = '0.0500'
$ws.Range("E43").Value = '  +3.97%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '0.997'
$ws.Range("E44").Value = '  -0.28%  '
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = '3.26'
$ws.Range("E45").Value = '  +0.96%  '
$ws.Range("E46").Value = '  +12.17%  '
$ws.Range("D47").Value = '3.57'
$ws.Range("E47").Value = '  +4.52%  '
$ws.Range("D48").Value = '2.70'
$ws.Range("D49").Value = '3.07'
$ws.Range("E49").Value = '  +2.56%  '
$ws.Range("D50").Value = '9.15'
$ws.Range("E50").Value = '  +7.30%  '
$ws.Range("D51").Value = '3.33'
$ws.Range("E51").Value = '  +2.07%  '
